$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly records for rows 2, 3, 5 and 6 were re-dated and had their
# Volumen / Precio / Origen data rotated between them, following the cycle:
#   new row2 <- old row5
#   new row5 <- old row3
#   new row3 <- old row6
#   new row6 <- old row2
# (row3/row6 keep their original "Origen" (column R) value, since it is the
# same text - "Región de La Araucanía" - before and after).
#
# Capture every source value first (using Value2 to avoid date/locale
# coercion on read), then perform all the writes, so the rotation is applied
# atomically and correctly regardless of evaluation order.

$D2 = $ws.Range("D2").Value2
$M2 = $ws.Range("M2").Value2
$N2 = $ws.Range("N2").Value2
$O2 = $ws.Range("O2").Value2
$P2 = $ws.Range("P2").Value2
$R2 = $ws.Range("R2").Value2
$S2 = $ws.Range("S2").Value2

$D3 = $ws.Range("D3").Value2
$M3 = $ws.Range("M3").Value2
$N3 = $ws.Range("N3").Value2
$O3 = $ws.Range("O3").Value2
$P3 = $ws.Range("P3").Value2
$R3 = $ws.Range("R3").Value2
$S3 = $ws.Range("S3").Value2

$D5 = $ws.Range("D5").Value2
$M5 = $ws.Range("M5").Value2
$N5 = $ws.Range("N5").Value2
$O5 = $ws.Range("O5").Value2
$P5 = $ws.Range("P5").Value2
$R5 = $ws.Range("R5").Value2
$S5 = $ws.Range("S5").Value2

$D6 = $ws.Range("D6").Value2
$M6 = $ws.Range("M6").Value2
$N6 = $ws.Range("N6").Value2
$O6 = $ws.Range("O6").Value2
$P6 = $ws.Range("P6").Value2
$S6 = $ws.Range("S6").Value2

# Row 2 <- old Row 5
$ws.Range("D2").Value = $D5
$ws.Range("M2").Value = $M5
$ws.Range("N2").Value = $N5
$ws.Range("O2").Value = $O5
$ws.Range("P2").Value = $P5
$ws.Range("R2").Value = $R5
$ws.Range("S2").Value = $S5

# Row 5 <- old Row 3 (Origen ends up "Región de La Araucanía", same text as
# row 3's own unchanged Origen)
$ws.Range("D5").Value = $D3
$ws.Range("M5").Value = $M3
$ws.Range("N5").Value = $N3
$ws.Range("O5").Value = $O3
$ws.Range("P5").Value = $P3
$ws.Range("R5").Value = $R3
$ws.Range("S5").Value = $S3

# Row 3 <- old Row 6
$ws.Range("D3").Value = $D6
$ws.Range("M3").Value = $M6
$ws.Range("N3").Value = $N6
$ws.Range("O3").Value = $O6
$ws.Range("P3").Value = $P6
$ws.Range("S3").Value = $S6

# Row 6 <- old Row 2
$ws.Range("D6").Value = $D2
$ws.Range("M6").Value = $M2
$ws.Range("N6").Value = $N2
$ws.Range("O6").Value = $O2
$ws.Range("P6").Value = $P2
$ws.Range("S6").Value = $S2
